$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H13").Value = 66670.664
$ws.Range("J13").Value = 66670.664
$ws.Range("L13").Value = 66670.664
$ws.Range("N13").Value = -67008.664

$ws.Range("H53").Value = 258.875
$ws.Range("I53").Value = 199.2
$ws.Range("J53").Value = 358.33334
$ws.Range("K53").Value = 199.2
$ws.Range("L53").Value = 358.33334
$ws.Range("M53").Value = 437.8
$ws.Range("N53").Value = -1632.33334

$ws.Range("H125").Value = 17827.715
$ws.Range("I125").Value = 30636
$ws.Range("J125").Value = 750
$ws.Range("K125").Value = 275724
$ws.Range("L125").Value = 6750
$ws.Range("M125").Value = -273264
$ws.Range("N125").Value = -11670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3219.2126
$ws.Range("I132").Value = 2657.9678
$ws.Range("J132").Value = 4306.625
$ws.Range("K132").Value = 7973.903399999999
$ws.Range("L132").Value = 12919.875
$ws.Range("M132").Value = -5443.903399999999
$ws.Range("N132").Value = -17979.875

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H45").Value = 1839.7222
$ws.Range("I45").Value = 1508.7059
$ws.Range("J45").Value = 2135.8948
$ws.Range("K45").Value = 1508.7059
$ws.Range("L45").Value = 2135.8948
$ws.Range("M45").Value = -1131.7059
$ws.Range("N45").Value = -2889.8948

$ws.Range("H95").Value = 29245.143
$ws.Range("J95").Value = 29245.143
$ws.Range("L95").Value = 29245.143
$ws.Range("N95").Value = -34737.143

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 1204.2354
$ws.Range("I122").Value = 815.5714
$ws.Range("K122").Value = 2446.7142
$ws.Range("M122").Value = 3.285799999999654

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 1600
$ws.Range("I23").Value = 400
$ws.Range("J23").Value = 2000
$ws.Range("K23").Value = 400
$ws.Range("L23").Value = 2000
$ws.Range("M23").Value = -160
$ws.Range("N23").Value = -2480

$ws.Range("H26").Value = 3000
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2713

$ws.Range("H27").Value = 1600
$ws.Range("I27").Value = 400
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 400
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -208
$ws.Range("N27").Value = -2384

$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H44").Value = 7400
$ws.Range("J44").Value = 7400
$ws.Range("L44").Value = 7400
$ws.Range("N44").Value = -8284

$ws.Range("H45").Value = 15000
$ws.Range("J45").Value = 15000
$ws.Range("L45").Value = 15000
$ws.Range("N45").Value = -16186

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H57").Value = 11000
$ws.Range("J57").Value = 11000
$ws.Range("L57").Value = 11000
$ws.Range("N57").Value = -12120

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 34488890
$ws.Range("I132").Value = 50007420
$ws.Range("J132").Value = 3256
$ws.Range("K132").Value = 150022260
$ws.Range("L132").Value = 9768
$ws.Range("M132").Value = -150019730
$ws.Range("N132").Value = -14828

$ws.Range("H22").Value = 223773.89
$ws.Range("I22").Value = 666757
$ws.Range("J22").Value = 2282.3333
$ws.Range("K22").Value = 2000271
$ws.Range("L22").Value = 6846.999899999999
$ws.Range("M22").Value = -2000102
$ws.Range("N22").Value = -7184.999899999999

$ws.Range("H27").Value = 223773.89
$ws.Range("I27").Value = 666757
$ws.Range("J27").Value = 2282.3333
$ws.Range("K27").Value = 2000271
$ws.Range("L27").Value = 6846.999899999999
$ws.Range("M27").Value = -2000169
$ws.Range("N27").Value = -7050.999899999999

$ws.Range("H113").Value = 14309.223
$ws.Range("I113").Value = 452.6
$ws.Range("J113").Value = 17458.455
$ws.Range("K113").Value = 1357.8
$ws.Range("L113").Value = 52375.36500000001
$ws.Range("M113").Value = 812.1999999999998
$ws.Range("N113").Value = -56715.36500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 1749.75
$ws.Range("I4").Value = 994.5
$ws.Range("J4").Value = 2505
$ws.Range("K4").Value = 994.5
$ws.Range("L4").Value = 2505
$ws.Range("M4").Value = -881.5
$ws.Range("N4").Value = -2731

$ws.Range("H22").Value = 3749.7812
$ws.Range("I22").Value = 557.2143
$ws.Range("J22").Value = 6232.8887
$ws.Range("K22").Value = 557.2143
$ws.Range("L22").Value = 6232.8887
$ws.Range("M22").Value = -262.2143
$ws.Range("N22").Value = -6822.8887

$ws.Range("H27").Value = 3749.7812
$ws.Range("I27").Value = 557.2143
$ws.Range("J27").Value = 6232.8887
$ws.Range("K27").Value = 557.2143
$ws.Range("L27").Value = 6232.8887
$ws.Range("M27").Value = -450.2143
$ws.Range("N27").Value = -6446.8887

$ws.Range("H28").Value = 1749.75
$ws.Range("I28").Value = 994.5
$ws.Range("J28").Value = 2505
$ws.Range("K28").Value = 994.5
$ws.Range("L28").Value = 2505
$ws.Range("M28").Value = -762.5
$ws.Range("N28").Value = -2969

$ws.Range("H32").Value = 406.5
$ws.Range("I32").Value = 406.5
$ws.Range("K32").Value = 406.5
$ws.Range("M32").Value = -89.5

$ws.Range("H33").Value = 9800
$ws.Range("J33").Value = 9800
$ws.Range("L33").Value = 9800
$ws.Range("N33").Value = -10380

$ws.Range("H37").Value = 1749.75
$ws.Range("I37").Value = 994.5
$ws.Range("J37").Value = 2505
$ws.Range("K37").Value = 994.5
$ws.Range("L37").Value = 2505
$ws.Range("M37").Value = -887.5
$ws.Range("N37").Value = -2719

$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4540

$ws.Range("H40").Value = 2403.4333
$ws.Range("I40").Value = 2071.25
$ws.Range("J40").Value = 2783.0715
$ws.Range("K40").Value = 2071.25
$ws.Range("L40").Value = 2783.0715
$ws.Range("M40").Value = -1935.25
$ws.Range("N40").Value = -3055.0715

$ws.Range("H46").Value = 1921.1
$ws.Range("I46").Value = 1014
$ws.Range("K46").Value = 1014
$ws.Range("M46").Value = -826

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H122").Value = 46961.348
$ws.Range("I122").Value = 93759.27
$ws.Range("J122").Value = 4063.25
$ws.Range("K122").Value = 281277.81
$ws.Range("L122").Value = 12189.75
$ws.Range("M122").Value = -278827.81
$ws.Range("N122").Value = -17089.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 11765867
$ws.Range("I122").Value = 15385395
$ws.Range("J122").Value = 2402.5
$ws.Range("K122").Value = 46156185
$ws.Range("L122").Value = 7207.5
$ws.Range("M122").Value = -46153735
$ws.Range("N122").Value = -12107.5
